# [LC-850] Release of LetsCo OS v1.3.0
# Rename the KPI short codes from GPn/BPn to GP0n/BP0n (2-digit suffix),
# and reset the sheet's scroll position / active selection back to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename KPI ids in column B (the "name" column of the rscKpi rows) ---
$ws.Range("B16").Value      = "GP01"   # was GP1
$ws.Range("B17:B18").Value  = "GP02"   # was GP2
$ws.Range("B19:B20").Value  = "GP03"   # was GP3
$ws.Range("B21:B25").Value  = "BP01"   # was BP1
$ws.Range("B26:B30").Value  = "BP02"   # was BP2
$ws.Range("B31:B60").Value  = "BP03"   # was BP3

# --- Reset the view: scroll back to A1 and make it the active selection ---
$ws.Range("A1").Select()
